$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 3:4 (pushes everything below down by 2, inheriting
# row 2's formatting for the new B3/B4 cells - same as Excel's native
# "insert row" behavior, which copies format from the row above).
$ws.Rows("3:4").Insert()

# Row 2: update the title text (style unchanged - still wrap-text s=2).
$ws.Range("B2").Value = "Mighty Aries\n<size=36>and the</size>\nMovement of Energy"

# Row 3 (new): company / RENEGADEWARE - keeps the wrap-text style inherited
# from the insert (matches target s="2").
$ws.Range("A3").Value = "company"
$ws.Range("B3").Value = "RENEGADEWARE"

# Row 4 (new): credits_music / long credits string - target has NO style on
# either cell, so clear the inherited formatting.
$ws.Range("A4").Value = "credits_music"
$ws.Range("B4").Value = "Music and Sound by Winfield B. Carson V\n\nSinger Katya Hall"
$ws.Range("A4:B4").ClearFormats()

# Update the selection to match the recorded workbook state.
$ws.Range("B2").Select()
